# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The G column ("K" = strikeouts) previously held pitch-count-derived
# "Strike#" values. It has been regenerated to hold the actual strikeout
# totals (K) for each start. Update the literal values in column G,
# rows 2-37, to the newly computed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 5
    4  = 10
    5  = 5
    6  = 10
    7  = 5
    8  = 7
    9  = 8
    10 = 4
    11 = 6
    12 = 6
    13 = 9
    14 = 9
    15 = 7
    16 = 7
    17 = 15
    18 = 6
    19 = 5
    20 = 4
    21 = 13
    22 = 11
    23 = 8
    24 = 4
    25 = 10
    26 = 9
    27 = 4
    28 = 13
    29 = 11
    30 = 5
    31 = 11
    32 = 7
    33 = 7
    34 = 6
    35 = 7
    36 = 2
    37 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
